$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Style donor for "text placeholder" cells (style index 14) ---
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null

# --- Set text placeholder values ---
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"

# --- Numeric cell updates ---
$ws.Range("D15").Value = 2
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 2
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J15").Value = 2
$ws.Range("J15").NumberFormat = '#,##0'
$ws.Range("K15").Value = -100
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -22.222222222222
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = 3.225806451612
$ws.Range("L16").Value = 113.333333333333
$ws.Range("M16").Value = 14.285714285714
$ws.Range("N16").Value = -81.609195402298
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -75
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -60
$ws.Range("I17").Value = 26
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -44.680851063829
$ws.Range("L17").Value = 23.809523809523
$ws.Range("M17").Value = 73.333333333333
$ws.Range("N17").Value = -56.666666666666
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -42.307692307692
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = -26.315789473684
$ws.Range("L18").Value = -12.5
$ws.Range("M18").Value = 14.285714285714
$ws.Range("N18").Value = -82.716049382716
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 90.909090909090
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 9.090909090909
$ws.Range("I19").Value = 167
$ws.Range("J19").Value = 133
$ws.Range("K19").Value = 25.563909774436
$ws.Range("L19").Value = 60.576923076923
$ws.Range("M19").Value = 103.658536585366
$ws.Range("N19").Value = 70.408163265306
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 12.5
$ws.Range("I20").Value = 37
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = 37.037037037037
$ws.Range("L20").Value = 105.555555555556
$ws.Range("M20").Value = 23.333333333333
$ws.Range("N20").Value = -83.482142857142
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -6.060606060606
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = -18.269230769230
$ws.Range("I21").Value = 318
$ws.Range("J21").Value = 316
$ws.Range("K21").Value = 0.632911392405
$ws.Range("L21").Value = 41.964285714285
$ws.Range("M21").Value = 55.121951219512
$ws.Range("N21").Value = -63.945578231292
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 10
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 25
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = 66.666666666666
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = -32.142857142857
$ws.Range("I24").Value = 230
$ws.Range("J24").Value = 270
$ws.Range("K24").Value = -14.814814814814
$ws.Range("L24").Value = 22.994652406417
$ws.Range("M24").Value = 84
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("I25").Value = 57
$ws.Range("J25").Value = 69
$ws.Range("K25").Value = -17.391304347826
$ws.Range("L25").Value = 7.547169811320
$ws.Range("M25").Value = 14
$ws.Range("D26").Value = 2
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G26").Value = 2
$ws.Range("G26").NumberFormat = '#,##0'
$ws.Range("H26").Value = -100
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J26").Value = 3
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 166.666666666667
$ws.Range("L27").Value = 77.777777777777
